$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as row 370, pushing the
# existing rows 370-398 down to 371-399 (dimension grows from R398 to R399).
$ws.Rows.Item(370).Insert()

$ws.Range("A370").Value = 3
$ws.Range("B370").Value = "Femacal de La Calera"
$ws.Range("C370").Value = "Coquimbo"
$ws.Range("D370").Value = 44783
$ws.Range("E370").Value = 5
$ws.Range("F370").Value = 100114013
$ws.Range("G370").Value = "Zanahoria"
$ws.Range("H370").Value = "Sin especificar"
$ws.Range("I370").Value = "Primera"
$ws.Range("J370").Value = 120
$ws.Range("K370").Value = 10000
$ws.Range("L370").Value = 10000
$ws.Range("M370").Value = 10000
$ws.Range("N370").Value = "`$/saco 20 kilos"
$ws.Range("O370").Value = "Provincia de Quillota"
$ws.Range("P370").Value = 500
$ws.Range("Q370").Value = 20
$ws.Range("R370").Value = "Hortaliza"
